$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (data reshuffled across the weekly logic) ---
# Row 2
$ws.Range("D2").Value = 44210
$ws.Range("K2").Value = 8000
$ws.Range("L2").Value = 9000
$ws.Range("M2").Value = 8417
$ws.Range("P2").Value = 140
# Row 3
$ws.Range("D3").Value = 44630
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 9000
$ws.Range("L3").Value = 9500
$ws.Range("M3").Value = 9250
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 154
# Row 4
$ws.Range("D4").Value = 44224
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 8500
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 8719
$ws.Range("P4").Value = 145
# Row 5
$ws.Range("D5").Value = 44216
$ws.Range("J5").Value = 55
$ws.Range("K5").Value = 9500
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = 9773
$ws.Range("P5").Value = 163
# Row 6
$ws.Range("D6").Value = 44690
$ws.Range("K6").Value = 10000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 10000
$ws.Range("O6").Value = "Región de Arica y Parinacota"
$ws.Range("P6").Value = 167
# Row 7
$ws.Range("D7").Value = 44271
$ws.Range("J7").Value = 55
$ws.Range("K7").Value = 9000
$ws.Range("L7").Value = 9500
$ws.Range("M7").Value = 9227
$ws.Range("P7").Value = 154
# Row 8
$ws.Range("D8").Value = 44253
$ws.Range("J8").Value = 95
$ws.Range("K8").Value = 9500
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = 9658
$ws.Range("P8").Value = 161
# Row 9
$ws.Range("D9").Value = 44204
$ws.Range("J9").Value = 45
$ws.Range("K9").Value = 9500
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = 9722
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 162
# Row 10
$ws.Range("D10").Value = 44208
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 7000
$ws.Range("M10").Value = 7350
$ws.Range("O10").Value = "Región del Maule"
$ws.Range("P10").Value = 122
# Row 11
$ws.Range("D11").Value = 44218
$ws.Range("J11").Value = 65
$ws.Range("M11").Value = 9615
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 160
# Row 12
$ws.Range("D12").Value = 44624
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 11000
$ws.Range("M12").Value = 10500
$ws.Range("P12").Value = 175
# Row 14
$ws.Range("D14").Value = 44617
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 11000
$ws.Range("M14").Value = 10500
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 175
# Row 15
$ws.Range("D15").Value = 44615
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 11000
$ws.Range("L15").Value = 12000
$ws.Range("M15").Value = 11500
$ws.Range("P15").Value = 192
# Row 16
$ws.Range("D16").Value = 44264
$ws.Range("J16").Value = 43
$ws.Range("K16").Value = 8500
$ws.Range("L16").Value = 9000
$ws.Range("M16").Value = 8709
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 145
# Row 17
$ws.Range("D17").Value = 44596
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 12000
$ws.Range("L17").Value = 13000
$ws.Range("M17").Value = 12500
$ws.Range("P17").Value = 208
# Row 18
$ws.Range("D18").Value = 44687
$ws.Range("K18").Value = 9000
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = 9500
$ws.Range("P18").Value = 158
# Row 19
$ws.Range("D19").Value = 44259
$ws.Range("J19").Value = 70
$ws.Range("K19").Value = 9000
$ws.Range("L19").Value = 9500
$ws.Range("M19").Value = 9214
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 154
# Row 20
$ws.Range("D20").Value = 44266
$ws.Range("J20").Value = 60
$ws.Range("K20").Value = 9000
$ws.Range("L20").Value = 9500
$ws.Range("M20").Value = 9208
$ws.Range("P20").Value = 153
# Row 21
$ws.Range("D21").Value = 44162
$ws.Range("J21").Value = 43
$ws.Range("K21").Value = 8000
$ws.Range("L21").Value = 8500
$ws.Range("M21").Value = 8209
$ws.Range("P21").Value = 137
# Row 23
$ws.Range("D23").Value = 44698
$ws.Range("J23").Value = 60
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = 10000
$ws.Range("P23").Value = 167
# Row 24
$ws.Range("D24").Value = 44627
# Row 25
$ws.Range("D25").Value = 44694
$ws.Range("J25").Value = 80
$ws.Range("K25").Value = 10000
$ws.Range("L25").Value = 10000
$ws.Range("M25").Value = 10000
$ws.Range("O25").Value = "Región de Arica y Parinacota"
$ws.Range("P25").Value = 167
# Row 26
$ws.Range("D26").Value = 44160
$ws.Range("J26").Value = 90
$ws.Range("K26").Value = 7500
$ws.Range("L26").Value = 8000
$ws.Range("M26").Value = 7667
$ws.Range("P26").Value = 128
# Row 27
$ws.Range("D27").Value = 44610
$ws.Range("K27").Value = 11000
$ws.Range("L27").Value = 12000
$ws.Range("M27").Value = 11500
$ws.Range("O27").Value = "Región Metropolitana"
$ws.Range("P27").Value = 192
# Row 28
$ws.Range("D28").Value = 44671
$ws.Range("J28").Value = 160
$ws.Range("K28").Value = 6000
$ws.Range("L28").Value = 7000
$ws.Range("M28").Value = 6500
$ws.Range("O28").Value = "Región de Arica y Parinacota"
$ws.Range("P28").Value = 108
# Row 29
$ws.Range("D29").Value = 44159
$ws.Range("J29").Value = 35
$ws.Range("K29").Value = 7500
$ws.Range("L29").Value = 8000
$ws.Range("M29").Value = 7714
$ws.Range("P29").Value = 129
# Row 30
$ws.Range("D30").Value = 44594
$ws.Range("J30").Value = 80
$ws.Range("K30").Value = 12000
$ws.Range("L30").Value = 13000
$ws.Range("M30").Value = 12500
$ws.Range("O30").Value = "Región de Arica y Parinacota"
$ws.Range("P30").Value = 208
# Row 31
$ws.Range("D31").Value = 44600
$ws.Range("J31").Value = 60
$ws.Range("K31").Value = 12000
$ws.Range("L31").Value = 13000
$ws.Range("M31").Value = 12500
$ws.Range("O31").Value = "Región de Arica y Parinacota"
$ws.Range("P31").Value = 208

# --- Add new row 32 ---
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C32").Value = "Ñuble"
$ws.Range("D32").Value = 44692
$ws.Range("D32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = 100112001
$ws.Range("G32").Value = "Berenjena"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 60
$ws.Range("K32").Value = 10000
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = 10000
$ws.Range("N32").Value = "`$/caja 60 unidades"
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 167
$ws.Range("Q32").Value = 60
$ws.Range("R32").Value = "Hortaliza"
